# Auto-generated edit script applying the OOXML diff cell-by-cell.
# Updates numeric values in the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# to match the values captured by the scheduled price-refresh job.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 981.41095
$ws.Range("I15").Value = 981.41095
$ws.Range("K15").Value = 2944.23285
$ws.Range("M15").Value = -2775.23285
$ws.Range("H100").Value = 8685.15
$ws.Range("I100").Value = 3177.7144
$ws.Range("J100").Value = 11650.692
$ws.Range("K100").Value = 3177.7144
$ws.Range("L100").Value = 11650.692
$ws.Range("M100").Value = -2636.7144
$ws.Range("N100").Value = -12732.692
$ws.Range("H107").Value = 1613.0834
$ws.Range("I107").Value = 1143.6
$ws.Range("K107").Value = 1143.6
$ws.Range("M107").Value = 776.4000000000001
$ws.Range("H121").Value = 2467.6
$ws.Range("J121").Value = 2467.6
$ws.Range("L121").Value = 7402.799999999999
$ws.Range("N121").Value = -10896.8
$ws.Range("H137").Value = 4516
$ws.Range("I137").Value = 5829.625
$ws.Range("K137").Value = 17488.875
$ws.Range("M137").Value = -14938.875
$ws.Range("H138").Value = 1694.63
$ws.Range("J138").Value = 1945.2764
$ws.Range("L138").Value = 5835.8292
$ws.Range("N138").Value = -16115.8292

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3202.5
$ws.Range("I45").Value = 2548.7273
$ws.Range("K45").Value = 2548.7273
$ws.Range("M45").Value = -2171.7273
$ws.Range("H61").Value = 1864495.4
$ws.Range("I61").Value = 4179073
$ws.Range("K61").Value = 4179073
$ws.Range("M61").Value = -4178861
$ws.Range("H107").Value = 93554.336
$ws.Range("J107").Value = 93554.336
$ws.Range("L107").Value = 93554.336
$ws.Range("N107").Value = -101234.336
$ws.Range("H132").Value = 2960901.8
$ws.Range("I132").Value = 4050875.2
$ws.Range("J132").Value = 2402.7144
$ws.Range("K132").Value = 12152625.6
$ws.Range("L132").Value = 7208.1432
$ws.Range("M132").Value = -12150095.6
$ws.Range("N132").Value = -12268.1432
$ws.Range("H136").Value = 1864495.4
$ws.Range("I136").Value = 4179073
$ws.Range("K136").Value = 12537219
$ws.Range("M136").Value = -12534669

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3577.6924
$ws.Range("I105").Value = 2630
$ws.Range("K105").Value = 2630
$ws.Range("M105").Value = -883

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2436.3447
$ws.Range("I31").Value = 2116.2173
$ws.Range("J31").Value = 3663.5
$ws.Range("K31").Value = 2116.2173
$ws.Range("L31").Value = 3663.5
$ws.Range("M31").Value = -1821.2173
$ws.Range("N31").Value = -4253.5
$ws.Range("H34").Value = 2436.3447
$ws.Range("I34").Value = 2116.2173
$ws.Range("J34").Value = 3663.5
$ws.Range("K34").Value = 2116.2173
$ws.Range("L34").Value = 3663.5
$ws.Range("M34").Value = -1914.2173
$ws.Range("N34").Value = -4067.5
$ws.Range("H58").Value = 4385.077
$ws.Range("J58").Value = 4503
$ws.Range("L58").Value = 4503
$ws.Range("N58").Value = -4909
$ws.Range("H94").Value = 2060.6428
$ws.Range("J94").Value = 2224.818
$ws.Range("L94").Value = 2224.818
$ws.Range("N94").Value = -3126.818
$ws.Range("H132").Value = 1742580.1
$ws.Range("I132").Value = 4447177.5
$ws.Range("K132").Value = 13341532.5
$ws.Range("M132").Value = -13339002.5
$ws.Range("H134").Value = 8423
$ws.Range("I134").Value = 7833
$ws.Range("K134").Value = 23499
$ws.Range("M134").Value = -20964
$ws.Range("H136").Value = 4385.077
$ws.Range("J136").Value = 4503
$ws.Range("L136").Value = 13509
$ws.Range("N136").Value = -18609

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 151.2963
$ws.Range("I2").Value = 173.1875
$ws.Range("J2").Value = 119.454544
$ws.Range("K2").Value = 1039.125
$ws.Range("L2").Value = 716.727264
$ws.Range("M2").Value = -926.125
$ws.Range("N2").Value = -942.727264
$ws.Range("H17").Value = 2217.2856
$ws.Range("J17").Value = 2086.8333
$ws.Range("L17").Value = 6260.499899999999
$ws.Range("N17").Value = -6598.499899999999
$ws.Range("H98").Value = 2966.3333
$ws.Range("I98").Value = 800
$ws.Range("J98").Value = 3399.6
$ws.Range("K98").Value = 2400
$ws.Range("L98").Value = 10198.8
$ws.Range("M98").Value = -902
$ws.Range("N98").Value = -13194.8
$ws.Range("H137").Value = 5325.3076
$ws.Range("J137").Value = 4249.875
$ws.Range("L137").Value = 12749.625
$ws.Range("N137").Value = -22949.625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 17500
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 25000
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 25000
$ws.Range("M52").Value = -9741
$ws.Range("N52").Value = -25518
$ws.Range("H97").Value = 1756.8077
$ws.Range("I97").Value = 881.2
$ws.Range("J97").Value = 2950.818
$ws.Range("K97").Value = 881.2
$ws.Range("L97").Value = 2950.818
$ws.Range("M97").Value = -385.2
$ws.Range("N97").Value = -3942.818
$ws.Range("H113").Value = 3008.4285
$ws.Range("I113").Value = 3008.4285
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3008.4285
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -838.4285
$ws.Range("N113").ClearContents()
$ws.Range("H123").Value = 44133
$ws.Range("J123").Value = 44133
$ws.Range("L123").Value = 44133
$ws.Range("N123").Value = -49033
$ws.Range("H124").Value = 106666.664
$ws.Range("J124").Value = 106666.664
$ws.Range("L124").Value = 106666.664
$ws.Range("N124").Value = -116486.664

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 984.9761999999999
$ws.Range("I16").Value = 524.5185
$ws.Range("J16").Value = 1813.8
$ws.Range("K16").Value = 524.5185
$ws.Range("L16").Value = 1813.8
$ws.Range("M16").Value = -354.5185
$ws.Range("N16").Value = -2153.8
$ws.Range("H22").Value = 2407
$ws.Range("I22").Value = 451.33334
$ws.Range("K22").Value = 451.33334
$ws.Range("M22").Value = -156.33334
$ws.Range("H27").Value = 2407
$ws.Range("I27").Value = 451.33334
$ws.Range("K27").Value = 451.33334
$ws.Range("M27").Value = -344.33334
$ws.Range("H100").Value = 3924.4666
$ws.Range("I100").Value = 3540.889
$ws.Range("J100").Value = 4499.8335
$ws.Range("K100").Value = 3540.889
$ws.Range("L100").Value = 4499.8335
$ws.Range("M100").Value = -2999.889
$ws.Range("N100").Value = -5581.8335
$ws.Range("H132").Value = 3161.2793
$ws.Range("I132").Value = 3144.9333
$ws.Range("K132").Value = 9434.7999
$ws.Range("M132").Value = -6904.7999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 96248.75
$ws.Range("J86").Value = 96248.75
$ws.Range("L86").Value = 96248.75
$ws.Range("N86").Value = -98494.75
$ws.Range("H89").Value = 96248.75
$ws.Range("J89").Value = 96248.75
$ws.Range("L89").Value = 481243.75
$ws.Range("N89").Value = -492475.75
$ws.Range("H96").Value = 3622
$ws.Range("J96").Value = 4699.75
$ws.Range("L96").Value = 4699.75
$ws.Range("N96").Value = -7445.75
$ws.Range("H100").Value = 875.0323
$ws.Range("I100").Value = 847.7037
$ws.Range("J100").Value = 1059.5
$ws.Range("K100").Value = 1695.4074
$ws.Range("L100").Value = 2119
$ws.Range("M100").Value = -1154.4074
$ws.Range("N100").Value = -3201
$ws.Range("H108").Value = 99998
$ws.Range("J108").Value = 99998
$ws.Range("L108").Value = 99998
$ws.Range("N108").Value = -107678
$ws.Range("H120").Value = 90000
$ws.Range("J120").Value = 90000
$ws.Range("L120").Value = 90000
$ws.Range("N120").Value = -99676
$ws.Range("H125").Value = 85857.14
$ws.Range("J125").Value = 85857.14
$ws.Range("L125").Value = 85857.14
$ws.Range("N125").Value = -95697.14
